# Auto-generated edit script: updates cryptocurrency price/volume data
# and fixes the ordering of two coin-pair rows (TRON/Polkadot, Filecoin/HuobiToken).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cell, $value) {
    # Force the cell to stay a literal text string, even when the value
    # looks like a number (e.g. "1.000", "311.48") so Excel does not
    # silently convert it to a numeric type and strip formatting.
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# Row 2
Set-TextCell $ws.Range("D2") "27.075.44"
Set-TextCell $ws.Range("E2") "  -2.48%  "

# Row 3
Set-TextCell $ws.Range("D3") "1.827.09"
Set-TextCell $ws.Range("E3") "  -1.30%  "

# Row 4
Set-TextCell $ws.Range("E4") "  -1.07%  "

# Row 5
Set-TextCell $ws.Range("D5") "311.48"
Set-TextCell $ws.Range("E5") "  -2.17%  "

# Row 6
Set-TextCell $ws.Range("D6") "1.001"
Set-TextCell $ws.Range("E6") "  -0.97%  "

# Row 7
Set-TextCell $ws.Range("E7") "  -1.88%  "

# Row 8
Set-TextCell $ws.Range("D8") "0.3675"
Set-TextCell $ws.Range("E8") "  -2.09%  "

# Row 9
Set-TextCell $ws.Range("D9") "0.07231"
Set-TextCell $ws.Range("E9") "  -1.63%  "

# Row 10
Set-TextCell $ws.Range("E10") "  -3.71%  "

# Row 11
Set-TextCell $ws.Range("D11") "20.76"
Set-TextCell $ws.Range("E11") "  -3.62%  "

# Row 12
Set-TextCell $ws.Range("D12") "1.820.18"
Set-TextCell $ws.Range("E12") "  -1.75%  "

# Row 13
Set-TextCell $ws.Range("D13") "6.662"
Set-TextCell $ws.Range("E13") "  -1.30%  "

# Row 14
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextCell $ws.Range("D14") "5.293"
Set-TextCell $ws.Range("E14") "  -2.78%  "

# Row 15
$ws.Range("B15").Value = "TRON"
$ws.Range("C15").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
Set-TextCell $ws.Range("D15") "0.07044"
Set-TextCell $ws.Range("E15") "  -1.16%  "

# Row 16
Set-TextCell $ws.Range("D16") "89.72"
Set-TextCell $ws.Range("E16") "  +0.83%  "

# Row 17
Set-TextCell $ws.Range("E17") "  -1.05%  "

# Row 18
Set-TextCell $ws.Range("D18") "0.000008755"
Set-TextCell $ws.Range("E18") "  -2.86%  "

# Row 19
Set-TextCell $ws.Range("E19") "  -0.96%  "

# Row 20
Set-TextCell $ws.Range("D20") "14.89"
Set-TextCell $ws.Range("E20") "  -3.61%  "

# Row 21
Set-TextCell $ws.Range("D21") "27.089.03"
Set-TextCell $ws.Range("E21") "  -2.49%  "

# Row 22
Set-TextCell $ws.Range("D22") "5.142"
Set-TextCell $ws.Range("E22") "  -1.42%  "

# Row 23
Set-TextCell $ws.Range("E23") "  -2.21%  "

# Row 24
Set-TextCell $ws.Range("D24") "2.046.79"
Set-TextCell $ws.Range("E24") "  -1.50%  "

# Row 25
Set-TextCell $ws.Range("D25") "1.984"
Set-TextCell $ws.Range("E25") "  +0.53%  "

# Row 26
Set-TextCell $ws.Range("D26") "151.66"
Set-TextCell $ws.Range("E26") "  -2.36%  "

# Row 27
Set-TextCell $ws.Range("E27") "  +3.69%  "

# Row 28
Set-TextCell $ws.Range("D28") "18.28"
Set-TextCell $ws.Range("E28") "  -1.88%  "

# Row 29
Set-TextCell $ws.Range("D29") "5.261"
Set-TextCell $ws.Range("E29") "  -1.95%  "

# Row 30
Set-TextCell $ws.Range("D30") "116.77"
Set-TextCell $ws.Range("E30") "  -1.75%  "

# Row 31
Set-TextCell $ws.Range("D31") "0.08702"
Set-TextCell $ws.Range("E31") "  -2.47%  "

# Row 32
Set-TextCell $ws.Range("E32") "  -4.00%  "

# Row 33
Set-TextCell $ws.Range("D33") "0.7377"
Set-TextCell $ws.Range("E33") "  -5.26%  "

# Row 34
$ws.Range("B34").Value = "HuobiToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-TextCell $ws.Range("D34") "2.899"
Set-TextCell $ws.Range("E34") "  -0.88%  "

# Row 35
$ws.Range("B35").Value = "Filecoin"
$ws.Range("C35").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextCell $ws.Range("D35") "4.426"
Set-TextCell $ws.Range("E35") "  -2.61%  "

# Row 36
Set-TextCell $ws.Range("D36") "1.000"
Set-TextCell $ws.Range("E36") "  -1.11%  "

# Row 37
Set-TextCell $ws.Range("D37") "1.092"
Set-TextCell $ws.Range("E37") "  -3.70%  "

# Row 38
Set-TextCell $ws.Range("D38") "0.05252"
Set-TextCell $ws.Range("E38") "  -1.88%  "

# Row 39
Set-TextCell $ws.Range("D39") "0.01944"
Set-TextCell $ws.Range("E39") "  -1.98%  "

# Row 40
Set-TextCell $ws.Range("D40") "7.319"
Set-TextCell $ws.Range("E40") "  +2.15%  "

# Row 41
Set-TextCell $ws.Range("D41") "2.880"
Set-TextCell $ws.Range("E41") "  -0.57%  "

# Row 42
Set-TextCell $ws.Range("D42") "0.1685"
Set-TextCell $ws.Range("E42") "  -0.72%  "

# Row 43
Set-TextCell $ws.Range("D43") "0.5071"
Set-TextCell $ws.Range("E43") "  -1.29%  "

# Row 44
Set-TextCell $ws.Range("D44") "8.551"
Set-TextCell $ws.Range("E44") "  -3.02%  "

# Row 45
Set-TextCell $ws.Range("D45") "10.52"
Set-TextCell $ws.Range("E45") "  -1.49%  "

# Row 46
Set-TextCell $ws.Range("D46") "1.954"
Set-TextCell $ws.Range("E46") "  +5.36%  "

# Row 47
Set-TextCell $ws.Range("D47") "106.07"
Set-TextCell $ws.Range("E47") "  -1.38%  "

# Row 48
Set-TextCell $ws.Range("D48") "0.4727"
Set-TextCell $ws.Range("E48") "  -0.79%  "

# Row 49
Set-TextCell $ws.Range("D49") "1.000"
Set-TextCell $ws.Range("E49") "  -1.12%  "

# Row 50
Set-TextCell $ws.Range("D50") "0.06331"
Set-TextCell $ws.Range("E50") "  -2.20%  "

# Row 51
Set-TextCell $ws.Range("D51") "1.656"
Set-TextCell $ws.Range("E51") "  -2.09%  "

